$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values: additional river counts (0 -> 1)
$ws.Range("G23").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 1

# Update the selection to match the new range (A1:O30), no active cell override
$ws.Range("A1:O30").Select()
